$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.957.05'
$ws.Range("E2").Value = '  -2.67%  '
$ws.Range("D3").Value = '1.859.51'
$ws.Range("E3").Value = '  -2.16%  '
$ws.Range("E4").Value = '  +0.23%  '
$ws.Range("D5").Value = "'305.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.91%  '
$ws.Range("E6").Value = '  +0.17%  '
$ws.Range("D7").Value = "'0.5044"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.64%  '
$ws.Range("D8").Value = "'0.3734"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.75%  '
$ws.Range("D9").Value = "'0.07128"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.49%  '
$ws.Range("E10").Value = '  -0.16%  '
$ws.Range("D11").Value = "'20.49"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.42%  '
$ws.Range("D12").Value = "'0.07547"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.82%  '
$ws.Range("D13").Value = '1.855.09'
$ws.Range("E13").Value = '  -4.35%  '
$ws.Range("D14").Value = "'5.289"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.74%  '
$ws.Range("D15").Value = "'88.91"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.86%  '
$ws.Range("E16").Value = '  +0.24%  '
$ws.Range("D17").Value = "'0.000008363"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.61%  '
$ws.Range("D18").Value = "'14.05"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.80%  '
$ws.Range("D19").Value = "'1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.13%  '
$ws.Range("D20").Value = '27.014.57'
$ws.Range("E20").Value = '  -2.57%  '
$ws.Range("D21").Value = "'5.047"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.59%  '
$ws.Range("D22").Value = '2.086.56'
$ws.Range("E22").Value = '  -4.12%  '
$ws.Range("E23").Value = '  -2.87%  '
$ws.Range("D24").Value = "'6.456"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'1.852"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.74%  '
$ws.Range("D26").Value = "'147.07"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.76%  '
$ws.Range("D27").Value = "'17.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.60%  '
$ws.Range("D28").Value = "'2.086"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.54%  '
$ws.Range("D29").Value = "'112.56"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.70%  '
$ws.Range("D30").Value = "'4.655"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.45%  '
$ws.Range("D31").Value = "'4.636"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.94%  '
$ws.Range("D32").Value = "'0.09020"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.28%  '
$ws.Range("D33").Value = "'0.05108"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.05%  '
$ws.Range("D34").Value = "'3.056"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.02%  '
$ws.Range("D35").Value = "'1.150"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.04%  '
$ws.Range("D36").Value = "'0.7275"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -7.10%  '
$ws.Range("D37").Value = "'0.02036"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.19%  '
$ws.Range("D38").Value = "'3.037"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.42%  '
$ws.Range("D39").Value = "'2.451"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -6.15%  '
$ws.Range("D40").Value = "'1.070"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.66%  '
$ws.Range("D41").Value = "'0.5312"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.13%  '
$ws.Range("D42").Value = "'6.579"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.88%  '
$ws.Range("D43").Value = "'115.07"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.59%  '
$ws.Range("D44").Value = "'8.277"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.81%  '
$ws.Range("D45").Value = "'0.1470"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.44%  '
$ws.Range("D46").Value = "'1.000"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.15%  '
$ws.Range("D47").Value = "'0.4603"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.23%  '
$ws.Range("D48").Value = "'10.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.42%  '
$ws.Range("D49").Value = "'1.558"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.13%  '
$ws.Range("E50").Value = '  -0.58%  '
$ws.Range("D51").Value = "'63.91"
$ws.Range("D51").Style = "Normal"
